# Generate Report for Handback
# Update the generated timestamps for the "4e117ccf-1c72-44d3-9990-b7b883e9169f" row
# on the Overview sheet, and the corresponding handoff/handback datetimes on the
# zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-09-03 06:48:10"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-09-03 06:48:03"
$zhcn.Range("K4").Value = "2016-09-03 06:48:31"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-09-03 06:48:10"
$dede.Range("K4").Value = "2016-09-03 06:48:39"
